$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> ECs)
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.014142
$ws.Range("H2").Value = 0.042426
$ws.Range("M2").Value = 6.436245333333333
$ws.Range("N2").Value = 19.308736
$ws.Range("O2").Value = 0.2367562936388591
$ws.Range("P2").Value = 0.2367562936388591
$ws.Range("Q2").Value = 0.09102138150399999
$ws.Range("R2").Value = 0.8191924335359999
$ws.Range("S2").Value = 0.2367562936388591
$ws.Range("T2").Value = 0.2367562936388591

# Row 3 (ECs -> FAPs)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.014142
$ws.Range("H3").Value = 0.042426
$ws.Range("O3").Value = 0.5508630013028089
$ws.Range("P3").Value = 0.550863001302809
$ws.Range("Q3").Value = 0.21178026834
$ws.Range("R3").Value = 1.90602241506
$ws.Range("S3").Value = 0.5508630013028089
$ws.Range("T3").Value = 0.550863001302809

# Row 4 (ECs -> MuSCs)
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.014142
$ws.Range("H4").Value = 0.042426
$ws.Range("O4").Value = 0.212380705058332
$ws.Range("P4").Value = 0.212380705058332
$ws.Range("Q4").Value = 0.08165014277799999
$ws.Range("R4").Value = 0.7348512850019999
$ws.Range("S4").Value = 0.212380705058332
$ws.Range("T4").Value = 0.212380705058332
